$p = $ppt.ActivePresentation
$p | Get-Member | Out-String | Write-Output
